$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Neutrino Enigma Unraveled: Unveiling the Ghost Particle's Secrets",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Exploring the Realm of Genetics: Unraveling the Secrets of Life", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Author name
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " Enrico Fermi",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Emily Watson", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Author e-mail address (three separate runs in the original document)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "enrico",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "emily", 2) | Out-Null

$d.Content.Find.Execute(
    "fermi@physicsinstitute",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "watson@schoolofbiology", 2) | Out-Null

$d.Content.Find.Execute(
    "org",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "edu", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) First body paragraph (three "paragraphs" separated by manual line
#    breaks) - replace the whole thing in one shot so the new wording
#    about genetics replaces the old wording about neutrinos.
# ---------------------------------------------------------------------------
$nl = [char]11
$bodyPara1 = $d.Paragraphs(5)
$rng1 = $bodyPara1.Range
$rng1.MoveEnd(1, -1) | Out-Null
$rng1.Text = (
    "In the intricate tapestry of life, genetics holds the key to understanding the symphony of inheritance." +
    " From the intricate dance of molecules to the vast canvas of biodiversity, this field unveils the enigmatic secrets of life." +
    " It delves into the blueprint of organisms, unravelling the mysteries of how traits and characteristics are passed down through generations." +
    "$nl${nl}At the heart of genetics lies the study of DNA, the molecule of life." +
    " DNA, with its double helix structure, acts as a blueprint for all living organisms." +
    " It contains the genetic instructions that determine an organism's traits and characteristics." +
    " The study of DNA and its interactions has led to groundbreaking discoveries in fields ranging from medicine to agriculture." +
    "$nl${nl}Genetics also explores the intricate world of gene expression." +
    " Genes, which are segments of DNA, contain the instructions for making proteins." +
    " Understanding how genes are expressed and regulated is essential for comprehending a wide range of biological processes, from development and growth to disease and evolution."
)

# ---------------------------------------------------------------------------
# 5) Summary paragraph - replace whole text (this also drops the
#    lastRenderedPageBreak that used to sit in front of it).
# ---------------------------------------------------------------------------
$bodyPara2 = $d.Paragraphs(7)
$rng2 = $bodyPara2.Range
$rng2.MoveEnd(1, -1) | Out-Null
$rng2.Text = (
    "Genetics holds the key to unraveling the mysteries of life, from the intricate dance of DNA to the wonders of biodiversity." +
    " Through the study of DNA and gene expression, this field uncovers the secrets of inheritance, variation, and biological processes." +
    " Genetics has revolutionized our understanding of life and continues to drive advancements in medicine, agriculture, and biotechnology."
)

# ---------------------------------------------------------------------------
# 6) Append a new, empty paragraph after the summary paragraph.
# ---------------------------------------------------------------------------
$d.Paragraphs(7).Range.InsertParagraphAfter() | Out-Null
